# hel-884: Correction du wording sur le taux de vétusté de constructions
# "Taux de vétusté des construction (en %)" -> "Taux de vétusté des constructions (en %)"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lisez-moi")

$oldText = "Taux de vétusté des construction (en %)"
$newText = "Taux de vétusté des constructions (en %)"

$cell = $ws.Cells.Find($oldText)
if ($cell) {
    $cell.Value = $newText
} else {
    # Fallback: the label historically lives in A18 of the "Lisez-moi" sheet.
    $ws.Range("A18").Value = $newText
}
